$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3182.0833
$ws.Range("J17").Value = 1581.3043
$ws.Range("L17").Value = 4743.9129
$ws.Range("N17").Value = -5079.9129
$ws.Range("H96").Value = 2266.2
$ws.Range("I96").Value = 5112.5
$ws.Range("J96").Value = 368.66666
$ws.Range("K96").Value = 15337.5
$ws.Range("L96").Value = 1105.99998
$ws.Range("M96").Value = -13964.5
$ws.Range("N96").Value = -3851.99998
$ws.Range("H99").Value = 304.16666
$ws.Range("I99").Value = 265
$ws.Range("K99").Value = 795
$ws.Range("M99").Value = 703
$ws.Range("H100").Value = 4834.357
$ws.Range("I100").Value = 4821.615
$ws.Range("K100").Value = 4821.615
$ws.Range("M100").Value = -4280.615
$ws.Range("H132").Value = 2294.2263
$ws.Range("I132").Value = 1675.8163
$ws.Range("J132").Value = 9869.75
$ws.Range("K132").Value = 5027.448899999999
$ws.Range("L132").Value = 29609.25
$ws.Range("M132").Value = -2497.448899999999
$ws.Range("N132").Value = -34669.25
$ws.Range("H138").Value = 2135.3777
$ws.Range("I138").Value = 1950.1666
$ws.Range("J138").Value = 2258.8518
$ws.Range("K138").Value = 5850.4998
$ws.Range("L138").Value = 6776.555399999999
$ws.Range("M138").Value = -710.4997999999996
$ws.Range("N138").Value = -17056.5554

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2518.2
$ws.Range("I45").Value = 2447.5
$ws.Range("J45").Value = 2565.3333
$ws.Range("K45").Value = 2447.5
$ws.Range("L45").Value = 2565.3333
$ws.Range("M45").Value = -2070.5
$ws.Range("N45").Value = -3319.3333
$ws.Range("H54").Value = 19000
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 19000
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 19000
$ws.Range("N54").Value = -20538
$ws.Range("M54").ClearContents()
$ws.Range("H74").Value = 24392124
$ws.Range("I74").Value = 26317518
$ws.Range("K74").Value = 26317518
$ws.Range("M74").Value = -26316644
$ws.Range("H77").Value = 24392124
$ws.Range("I77").Value = 26317518
$ws.Range("K77").Value = 131587590
$ws.Range("M77").Value = -131583222
$ws.Range("H97").Value = 1349.8667
$ws.Range("I97").Value = 1349.8667
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 1349.8667
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -853.8667
$ws.Range("N97").ClearContents()
$ws.Range("H102").Value = 2639682
$ws.Range("I102").Value = 2937814
$ws.Range("J102").Value = 254624.75
$ws.Range("K102").Value = 2937814
$ws.Range("L102").Value = 254624.75
$ws.Range("M102").Value = -2936192
$ws.Range("N102").Value = -257868.75
$ws.Range("H110").Value = 47621024
$ws.Range("I110").Value = 55557310
$ws.Range("J110").Value = 3299.3333
$ws.Range("K110").Value = 55557310
$ws.Range("L110").Value = 3299.3333
$ws.Range("M110").Value = -55555265
$ws.Range("N110").Value = -7389.3333
$ws.Range("H122").Value = 2184.45
$ws.Range("I122").Value = 1499.96
$ws.Range("K122").Value = 4499.88
$ws.Range("M122").Value = -2049.88
$ws.Range("H132").Value = 2358.158
$ws.Range("I132").Value = 1326.9333
$ws.Range("K132").Value = 3980.7999
$ws.Range("M132").Value = -1450.7999

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 9260696
$ws.Range("I94").Value = 10417684
$ws.Range("J94").Value = 4789
$ws.Range("K94").Value = 10417684
$ws.Range("L94").Value = 4789
$ws.Range("M94").Value = -10417233
$ws.Range("N94").Value = -5691

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 994.8889
$ws.Range("I22").Value = 520.8
$ws.Range("J22").Value = 1587.5
$ws.Range("K22").Value = 520.8
$ws.Range("L22").Value = 1587.5
$ws.Range("M22").Value = -170.8
$ws.Range("N22").Value = -2287.5
$ws.Range("H31").Value = 2727.1482
$ws.Range("I31").Value = 2143.5334
$ws.Range("J31").Value = 2951.6155
$ws.Range("K31").Value = 2143.5334
$ws.Range("L31").Value = 2951.6155
$ws.Range("M31").Value = -1848.5334
$ws.Range("N31").Value = -3541.6155
$ws.Range("H34").Value = 2727.1482
$ws.Range("I34").Value = 2143.5334
$ws.Range("J34").Value = 2951.6155
$ws.Range("K34").Value = 2143.5334
$ws.Range("L34").Value = 2951.6155
$ws.Range("M34").Value = -1941.5334
$ws.Range("N34").Value = -3355.6155
$ws.Range("H58").Value = 2874.182
$ws.Range("I58").Value = 1306.0454
$ws.Range("K58").Value = 1306.0454
$ws.Range("M58").Value = -1103.0454
$ws.Range("H122").Value = 2685.2144
$ws.Range("I122").Value = 3345.5557
$ws.Range("J122").Value = 1496.6
$ws.Range("K122").Value = 10036.6671
$ws.Range("L122").Value = 4489.799999999999
$ws.Range("M122").Value = -7586.667099999999
$ws.Range("N122").Value = -9389.799999999999
$ws.Range("H136").Value = 2874.182
$ws.Range("I136").Value = 1306.0454
$ws.Range("K136").Value = 3918.1362
$ws.Range("M136").Value = -1368.1362

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 40000504
$ws.Range("I122").Value = 640.6667
$ws.Range("J122").Value = 100000300
$ws.Range("K122").Value = 5766.0003
$ws.Range("L122").Value = 900002700
$ws.Range("M122").Value = -3316.0003
$ws.Range("N122").Value = -900007600
$ws.Range("H132").Value = 2571.0334
$ws.Range("I132").Value = 1748.9166
$ws.Range("J132").Value = 2776.5625
$ws.Range("K132").Value = 15740.2494
$ws.Range("L132").Value = 24989.0625
$ws.Range("M132").Value = -13210.2494
$ws.Range("N132").Value = -30049.0625

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("N35").ClearContents()
$ws.Range("H97").Value = 695.82355
$ws.Range("I97").Value = 695.82355
$ws.Range("K97").Value = 695.82355
$ws.Range("M97").Value = -199.82355
$ws.Range("H113").Value = 3449.2942
$ws.Range("I113").Value = 2613.8
$ws.Range("K113").Value = 2613.8
$ws.Range("M113").Value = -443.8000000000002

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 851.5238000000001
$ws.Range("I82").Value = 577.1818
$ws.Range("K82").Value = 577.1818
$ws.Range("M82").Value = -216.1818
$ws.Range("H85").Value = 851.5238000000001
$ws.Range("I85").Value = 577.1818
$ws.Range("K85").Value = 577.1818
$ws.Range("M85").Value = 670.8182
$ws.Range("H100").Value = 3123.25
$ws.Range("I100").Value = 3123.25
$ws.Range("K100").Value = 3123.25
$ws.Range("M100").Value = -2582.25
$ws.Range("H110").Value = 28459
$ws.Range("J110").Value = 28459
$ws.Range("L110").Value = 28459
$ws.Range("N110").Value = -36639
$ws.Range("H122").Value = 4427.933
$ws.Range("I122").Value = 3964.9
$ws.Range("K122").Value = 11894.7
$ws.Range("M122").Value = -9444.700000000001
$ws.Range("H132").Value = 34489148
$ws.Range("I132").Value = 45457900
$ws.Range("K132").Value = 136373700
$ws.Range("M132").Value = -136371170

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 28577714
$ws.Range("I81").Value = 500
$ws.Range("J81").Value = 40008600
$ws.Range("K81").Value = 1000
$ws.Range("L81").Value = 80017200
$ws.Range("M81").Value = 61
$ws.Range("N81").Value = -80019322
$ws.Range("H84").Value = 28577714
$ws.Range("I84").Value = 500
$ws.Range("J84").Value = 40008600
$ws.Range("K84").Value = 5000
$ws.Range("L84").Value = 400086000
$ws.Range("M84").Value = 304
$ws.Range("N84").Value = -400096608
$ws.Range("H100").Value = 27780148
$ws.Range("I100").Value = 31252554
$ws.Range("J100").Value = 900
$ws.Range("K100").Value = 62505108
$ws.Range("L100").Value = 1800
$ws.Range("M100").Value = -62504567
$ws.Range("N100").Value = -2882
$ws.Range("H122").Value = 1634.5217
$ws.Range("I122").Value = 1493.5667
$ws.Range("K122").Value = 4480.7001
$ws.Range("M122").Value = -2030.7001
$ws.Range("H126").Value = 1456.174
$ws.Range("I126").Value = 1088.1177
$ws.Range("J126").Value = 2499
$ws.Range("K126").Value = 3264.3531
$ws.Range("L126").Value = 7497
$ws.Range("M126").Value = -794.3531000000003
$ws.Range("N126").Value = -12437
$ws.Range("H132").Value = 4073.2036
$ws.Range("I132").Value = 4015.4888
$ws.Range("J132").Value = 4361.778
$ws.Range("K132").Value = 12046.4664
$ws.Range("L132").Value = 13085.334
$ws.Range("M132").Value = -9516.466400000001
$ws.Range("N132").Value = -18145.334
